$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.757764577865601
$ws.Range("B1").Value = 3.787448406219482
$ws.Range("C1").Value = 1.940144777297974
$ws.Range("D1").Value = 1.341444969177246
$ws.Range("E1").Value = 1.135859489440918
